# Update the "取得日時" (acquired timestamp) column for all data rows on the
# "ランサーズ" sheet from 2025-11-21 18:23:52 to 2025-11-21 18:30:41,
# reflecting a fresh scrape appended at that time (commit: "Append: 2025-11-21 18:30 JST").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-11-21 18:23:52"
$newTimestamp = "2025-11-21 18:30:41"

# Data starts at row 2 (row 1 is the header) and ends at the last used row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row()

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value() -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
